# Fruta / hortaliza, semanal
# Insert a new weekly record at row 21 of the "Tuna" sheet, pushing the
# existing rows (21-91) down to (22-92).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before row 21; Excel shifts rows 21-91 -> 22-92
# and copies row 21's formatting (incl. the date-formatted style on D)
# onto the new row.
$ws.Rows.Item(21).Insert()

# Populate the new row 21 with the new weekly observation.
$ws.Range("A21").Value = 10
$ws.Range("B21").Value = "Vega Modelo de Temuco"
$ws.Range("C21").Value = "La Araucanía"
$ws.Range("D21").Value = 45037
$ws.Range("E21").Value = 9
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100107
$ws.Range("H21").Value = "Otros"
$ws.Range("I21").Value = 100107011
$ws.Range("J21").Value = "Tuna"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 70
$ws.Range("N21").Value = 22000
$ws.Range("O21").Value = 22000
$ws.Range("P21").Value = 22000
$ws.Range("Q21").Value = "$/caja 16 kilos"
$ws.Range("R21").Value = "Provincia de Los Andes"
$ws.Range("S21").Value = 1375
$ws.Range("T21").Value = 16
